# Applies the README/docx stats fix for the Renaissance - JDK17 - ZGC
# gauss-mix / heap-8G benchmark table.
#
# Summary of the change:
#  - Rows 1-3 (first three data rows) become placeholder "0M" values.
#  - Ten brand-new single-value rows are inserted right after them,
#    carrying the real per-iteration stats that used to be crammed,
#    tab-separated, into the three rows near the end of the table.
#  - Those three old "summary" rows near the end are collapsed down to
#    a single plain value each (the tab-separated figures are dropped).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Sanity-check the row layout we expect to find before mutating anything.
# (Cell.Range.Text carries a trailing cell-mark, hence the -like wildcard.)
if ($t.Cell(1, 1).Range.Text -notlike "99.99*") {
    throw "Unexpected content in row 1 - aborting."
}
if ($t.Cell(2, 1).Range.Text -notlike "0.01*") {
    throw "Unexpected content in row 2 - aborting."
}
if ($t.Cell(3, 1).Range.Text -notlike "88*") {
    throw "Unexpected content in row 3 - aborting."
}

# --- 1. Insert the 10 new rows right before current row 4 ("0") -----------
# NOTE: Rows.Add(beforeRow) always inserts immediately in front of the row
# object you hand it. Re-using the very same anchor for every call would
# therefore build the block in *reverse*, so the anchor index is advanced
# by one after each insertion to keep the values in their intended order.
$newValues = @("104", "0.00003", "0.00012", "0.00006", "0.00002", "0.00004", "0.00012", "0.00012", "0.00509", "100.0")

$anchorIndex = 4
foreach ($val in $newValues) {
    $beforeRow = $t.Rows.Item($anchorIndex)
    $newRow = $t.Rows.Add($beforeRow)
    $newRow.Cells.Item(1).Range.Text = $val
    $anchorIndex = $anchorIndex + 1
}

# --- 2. Collapse the first three rows to "0M" ------------------------------
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# --- 3. Collapse the three tab-separated summary rows near the end --------
# After the 10-row insertion, everything from the old row 4 onward shifted
# down by 10, so the three rows that used to be 34/35/36 are now 44/45/46.
$t.Cell(44, 1).Range.Text = "99.99"
$t.Cell(45, 1).Range.Text = "0.01"
$t.Cell(46, 1).Range.Text = "88"

Write-Host "Final row count:" $t.Rows.Count
